# Applies updated label values (rank/color/spatial metrics) to rows 7-9
# of the aggregate metrics sheet, matching a rerun of the underlying plots.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 17
$ws.Range("L7").Value = 7
$ws.Range("M7").Value = 17
$ws.Range("O7").Value = 7
$ws.Range("P7").Value = 17
$ws.Range("R7").Value = 7
$ws.Range("S7").Value = 17
$ws.Range("U7").Value = 7
$ws.Range("V7").Value = 17
$ws.Range("X7").Value = 7
$ws.Range("Y7").Value = 17
$ws.Range("AA7").Value = 7
$ws.Range("AB7").Value = 17
$ws.Range("AD7").Value = 7
$ws.Range("AE7").Value = 16
$ws.Range("AG7").Value = 7
$ws.Range("AH7").Value = 17
$ws.Range("AJ7").Value = 7
$ws.Range("AK7").Value = 15
$ws.Range("AM7").Value = 7
$ws.Range("AN7").Value = 15
$ws.Range("AP7").Value = 119
$ws.Range("AQ7").Value = 239
$ws.Range("AS7").Value = 88
$ws.Range("AT7").Value = 239
$ws.Range("AV7").Value = 75
$ws.Range("AW7").Value = 240
$ws.Range("AY7").Value = 65
$ws.Range("AZ7").Value = 230
$ws.Range("BB7").Value = 56
$ws.Range("BE7").Value = 48
$ws.Range("BG7").Value = 0.2431951060327109
$ws.Range("BH7").Value = 0.5016769669062087
$ws.Range("BI7").Value = 1
$ws.Range("BJ7").Value = 0
$ws.Range("BK7").Value = 0.406294602583436
$ws.Range("BL7").Value = 0.9929078014184397
$ws.Range("BN7").Value = 0.1497718489638923
$ws.Range("BO7").Value = 0.3938661711008237
$ws.Range("BQ7").Value = 0.1542546444780535
$ws.Range("BR7").Value = 0.4435346124233666
$ws.Range("BS7").Value = [double]"4.585473220836391e-05"
$ws.Range("BT7").Value = 0.04866770861981455
$ws.Range("BU7").Value = 0.5152027027027027

# Row 8
$ws.Range("D8").Value = $false
$ws.Range("I8").Value = 6.5
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 6.5
$ws.Range("Q8").Value = 2
$ws.Range("S8").Value = 8
$ws.Range("T8").Value = 2
$ws.Range("V8").Value = 9
$ws.Range("W8").Value = 4
$ws.Range("X8").Value = 6
$ws.Range("Z8").Value = 3
$ws.Range("AA8").Value = 5
$ws.Range("AB8").Value = 13
$ws.Range("AC8").Value = 1
$ws.Range("AD8").Value = 4.5
$ws.Range("AE8").Value = 11
$ws.Range("AG8").Value = 7
$ws.Range("AH8").Value = 12
$ws.Range("AJ8").Value = 6
$ws.Range("AK8").Value = 11
$ws.Range("AL8").Value = 2
$ws.Range("AM8").Value = 6
$ws.Range("AN8").Value = 11
$ws.Range("AO8").Value = 6
$ws.Range("AP8").Value = 108
$ws.Range("AQ8").Value = 179
$ws.Range("AR8").Value = 6
$ws.Range("AS8").Value = 75.5
$ws.Range("AT8").Value = 134
$ws.Range("AU8").Value = 15
$ws.Range("AV8").Value = 62
$ws.Range("AW8").Value = 108
$ws.Range("AX8").Value = 0
$ws.Range("AY8").Value = 55
$ws.Range("AZ8").Value = 95
$ws.Range("BA8").Value = 0
$ws.Range("BB8").Value = 44.5
$ws.Range("BC8").Value = 79
$ws.Range("BD8").Value = 0
$ws.Range("BE8").Value = 34.5
$ws.Range("BF8").Value = 93
$ws.Range("BG8").Value = 0.4654117848285294
$ws.Range("BH8").Value = 0.5443347513181849
$ws.Range("BI8").Value = 0.9870761866773036
$ws.Range("BJ8").Value = 0.303225806451613
$ws.Range("BK8").Value = 0.4294966012707016
$ws.Range("BL8").Value = 0.7254102644507164
$ws.Range("BM8").Value = 0
$ws.Range("BN8").Value = 0.1578496260175382
$ws.Range("BO8").Value = 0.3441590921303435
$ws.Range("BP8").Value = 0
$ws.Range("BQ8").Value = 0.2107149380697658
$ws.Range("BR8").Value = 0.3409516951580344
$ws.Range("BS8").Value = [double]"2.858286171611502e-05"
$ws.Range("BT8").Value = 0.04341140999265811
$ws.Range("BU8").Value = 0.3021972884525479

# Row 9
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 14
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 6
$ws.Range("M9").Value = 12
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 6
$ws.Range("R9").Value = 6
$ws.Range("S9").Value = 14
$ws.Range("U9").Value = 6
$ws.Range("V9").Value = 13
$ws.Range("X9").Value = 6
$ws.Range("Y9").Value = 13
$ws.Range("AB9").Value = 15
$ws.Range("AD9").Value = 6
$ws.Range("AE9").Value = 17
$ws.Range("AH9").Value = 14
$ws.Range("AK9").Value = 17
$ws.Range("AN9").Value = 17
$ws.Range("AO9").Value = 1
$ws.Range("AP9").Value = 115.5
$ws.Range("AQ9").Value = 232
$ws.Range("AR9").Value = 1
$ws.Range("AS9").Value = 84
$ws.Range("AT9").Value = 230
$ws.Range("AU9").Value = 2
$ws.Range("AV9").Value = 70.5
$ws.Range("AW9").Value = 235
$ws.Range("AX9").Value = 0
$ws.Range("AY9").Value = 61
$ws.Range("AZ9").Value = 201
$ws.Range("BA9").Value = 0
$ws.Range("BB9").Value = 53
$ws.Range("BC9").Value = 203
$ws.Range("BD9").Value = 0
$ws.Range("BE9").Value = 45
$ws.Range("BF9").Value = 207
$ws.Range("BG9").Value = 0.2629369803476946
$ws.Range("BH9").Value = 0.5079363171785225
$ws.Range("BI9").Value = 1
$ws.Range("BJ9").Value = 0.1022727272727273
$ws.Range("BK9").Value = 0.4255890480611942
$ws.Range("BL9").Value = 0.885304659498208
$ws.Range("BM9").Value = 0
$ws.Range("BN9").Value = 0.152633952419462
$ws.Range("BO9").Value = 0.3668940743409324
$ws.Range("BP9").Value = 0
$ws.Range("BQ9").Value = 0.1675755632780588
$ws.Range("BR9").Value = 0.3851287578467805
$ws.Range("BS9").Value = [double]"4.625346901017576e-05"
$ws.Range("BT9").Value = 0.05646631042973063
$ws.Range("BU9").Value = 0.5655471289274107
